# Updated symbol list on Wed Feb 15 15:27:40 UTC 2023 with GitHub Actions
# Refresh the crypto "Price" (column D) and "Volume(1h)" (column E) figures
# for each coin row on Sheet1. Source values are text (not numbers/percentages
# as Excel types), so each cell is forced to Text format before the write and
# restored to the Normal cell style afterwards so only the displayed text
# changes and no stray number-format style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2'  = '302.10';       'E2'  = '1.06%'
    'D3'  = '43.31';        'E3'  = '5.44%'
    'D4'  = '5.088';        'E4'  = '1.02%'
    'D5'  = '0.07699';      'E5'  = '2.75%'
    'D6'  = '1.619';        'E6'  = '2.06%'
    'D7'  = '1.036';        'E7'  = '10.28%'
    'D8'  = '0.1250';       'E8'  = '3.19%'
    'D9'  = '0.1853';       'E9'  = '2.00%'
    'D10' = '0.09121';      'E10' = '3.91%'
    'D11' = '0.04168';      'E11' = '-1.82%'
    'D12' = '0.1047';       'E12' = '-0.26%'
    'D13' = '0.001281';     'E13' = '1.59%'
    'D14' = '0.005748';     'E14' = '-0.24%'
    'E15' = '1,908.40%'
    'D16' = '3.331';        'E16' = '-0.80%'
    'D17' = '4.415';        'E17' = '1.41%'
    'E18' = '-1.92%'
    'D19' = '0.3355';       'E19' = '1.44%'
    'D20' = '8.562';        'E20' = '8.31%'
    'D21' = '0.1367';       'E21' = '-0.37%'
    'D22' = '0.3196';       'E22' = '8.43%'
    'D23' = '0.04154';      'E23' = '3.71%'
    'D24' = '0.001285';     'E24' = '1.88%'
    'D25' = '0.004479';     'E25' = '15.56%'
    'E26' = '10.28%'
    'D38' = '0.02454';      'E38' = '1.24%'
    'D39' = '0.05267';      'E39' = '1.96%'
    'D40' = '0.005962';     'E40' = '-1.74%'
    'D41' = '0.007668';     'E41' = '-0.90%'
    'D42' = '0.1346';       'E42' = '1.59%'
    'D43' = '0.007351';     'E43' = '0.43%'
    'D44' = '0.007552';     'E44' = '5.47%'
    'D45' = '0.3012';       'E45' = '1.24%'
    'D46' = '0.00006706';   'E46' = '7.62%'
    'D47' = '0.00000000749';'E47' = '0.36%'
    'D48' = '0.04480';      'E48' = '-4.32%'
    'E49' = '0.62%'
    'D50' = '0.00002097';   'E50' = '0.36%'
    'D51' = '0.0001997';    'E51' = '0.36%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force Text storage so the literal digits/percent sign survive exactly
    # (otherwise Excel would silently coerce "302.10" -> 302.1, "1.06%" -> 0.0106).
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
    # Drop back to the Normal cell style so no explicit "Text" number-format
    # style index is left behind on the cell (matches original formatting).
    $cell.Style = 'Normal'
}
